$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.645.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -6.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.300.17"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.82%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -10.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "612.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.41%  "

$ws.Range("E7").Value = "  -11.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.368"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -13.20%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.880"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -15.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.300.44"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.187"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -11.03%  "

$ws.Range("E13").Value = "  -14.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "91.584.20"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.71"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -7.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.911.73"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.19%  "

$ws.Range("E17").Value = "  -9.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.296.38"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -14.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -12.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -10.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "475.02"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -8.76%  "

$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -11.95%  "

$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.419"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -19.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000174"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -13.42%  "

$ws.Range("E26").Value = "  -12.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "87.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -10.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -12.60%  "

$ws.Range("E29").Value = "  +0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("E32").Value = "  -10.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.126"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -12.60%  "

$ws.Range("E34").Value = "  -12.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.61"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -10.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.501"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -14.78%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "508.90"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.14"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.33"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -11.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.143"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -8.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.93"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.844"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.78%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.12%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.06"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "51.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0377"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -12.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -12.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.78%  "
